$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (233) down into the new rows (234:238)
$src = $ws.Range("A233:D233")
$dst = $ws.Range("A234:D238")
$src.Copy()
$dst.PasteSpecial(-4122)

# New data rows to append (column A is an Excel date serial, formatted as date/time)
$newRows = @(
    @{ Row = 234; A = 44308; B = 3; C = 31; D = 365.2645222104395 },
    @{ Row = 235; A = 44309; B = 2; C = 25; D = 294.5681630729351 },
    @{ Row = 236; A = 44310; B = 1; C = 19; D = 223.8718039354306 },
    @{ Row = 237; A = 44311; B = 0; C = 17; D = 200.3063508895958 },
    @{ Row = 238; A = 44312; B = 0; C = 17; D = 200.3063508895958 }
)

foreach ($r in $newRows) {
    $rowIndex = $r.Row
    $ws.Cells.Item($rowIndex, 1).Value2 = $r.A
    $ws.Cells.Item($rowIndex, 2).Value2 = $r.B
    $ws.Cells.Item($rowIndex, 3).Value2 = $r.C
    $ws.Cells.Item($rowIndex, 4).Value2 = $r.D
}
